$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.100.33'
$ws.Range("E2").Value = '  -1.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.071.48'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.87'
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.90'
$ws.Range("E6").Value = '  +4.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.546'
$ws.Range("E8").Value = '  +3.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.077.63'
$ws.Range("E9").Value = '  -1.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -3.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.88'
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("E12").Value = '  +0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("E13").Value = '  -1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.32'
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("E15").Value = '  -1.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.579.57'
$ws.Range("E16").Value = '  -2.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.23'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.212.35'
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.073.74'
$ws.Range("E19").Value = '  -2.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.85'
$ws.Range("E20").Value = '  +2.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.70'
$ws.Range("E21").Value = '  +2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.720'
$ws.Range("E22").Value = '  -1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.56'
$ws.Range("E23").Value = '  +0.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  +3.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.18'
$ws.Range("E25").Value = '  +1.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.51'
$ws.Range("E26").Value = '  +0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  +3.39%  '

$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.68'
$ws.Range("E30").Value = '  -1.12%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  -2.08%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.995'
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  +4.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.32'
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0855'
$ws.Range("E35").Value = '  +2.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  -0.72%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.39'
$ws.Range("E37").Value = '  +4.61%  '

$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.15'
$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("E39").Value = '  -3.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.34'
$ws.Range("E40").Value = '  +1.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.39'
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '447.37'
$ws.Range("E42").Value = '  -2.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  -2.31%  '

$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.60'
$ws.Range("E44").Value = '  +2.60%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0365'
$ws.Range("E45").Value = '  -2.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.816.39'
$ws.Range("E46").Value = '  -3.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.111'
$ws.Range("E47").Value = '  +2.22%  '

$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.35'
$ws.Range("E49").Value = '  +5.05%  '

$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.27'
$ws.Range("E51").Value = '  +0.65%  '
